$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.004240061173064191
$ws.Range("D2").Value = 0.001023183206871892
$ws.Range("E2").Value = 0.04002618873340502
$ws.Range("F2").Value = 0.6715747400855463
$ws.Range("G2").Value = 0.002357098138867389
$ws.Range("I2").Value = 0.5002110225797054
$ws.Range("M2").Value = 0.9687595036895544
$ws.Range("N2").Value = 1.748505083762183
$ws.Range("O2").Value = 2.190178209650639
$ws.Range("C3").Value = 0.003698286982462662
$ws.Range("D3").Value = 0.001052514461458376
$ws.Range("E3").Value = 0.04258201232814418
$ws.Range("F3").Value = 0.6338536915487367
$ws.Range("G3").Value = 0.002361219472908731
$ws.Range("I3").Value = 0.4720757699876614
$ws.Range("M3").Value = 0.8536493934495439
$ws.Range("N3").Value = 1.603114105708613
$ws.Range("O3").Value = 2.072273628073162
$ws.Range("C4").Value = 0.003365212243334526
$ws.Range("D4").Value = 0.001071733232506755
$ws.Range("E4").Value = 0.04427033815240478
$ws.Range("F4").Value = 0.6111867210869946
$ws.Range("G4").Value = 0.002363885648601975
$ws.Range("I4").Value = 0.4551717564147779
$ws.Range("M4").Value = 0.7828418212515231
$ws.Range("N4").Value = 1.514014011049994
$ws.Range("O4").Value = 2.001542061349653
$ws.Range("C5").Value = 0.003229380730392251
$ws.Range("D5").Value = 0.001079865354451925
$ws.Range("E5").Value = 0.04498802544245351
$ws.Range("F5").Value = 0.6020730153911416
$ws.Range("G5").Value = 0.002365006362123455
$ws.Range("I5").Value = 0.4483758417128527
$ws.Range("M5").Value = 0.7539560113548305
$ws.Range("N5").Value = 1.477752157526254
$ws.Range("O5").Value = 1.973133107656452
$ws.Range("C6").Value = 0.003206820084134421
$ws.Range("D6").Value = 0.001081233716898389
$ws.Range("E6").Value = 0.04510898214972858
$ws.Range("F6").Value = 0.6005671096304184
$ws.Range("G6").Value = 0.002365194526189596
$ws.Range("I6").Value = 0.4472529576121502
$ws.Range("M6").Value = 0.7491577053726246
$ws.Range("N6").Value = 1.471733879279242
$ws.Range("O6").Value = 1.96844077623274
$ws.Range("C7").Value = 0.003363380772317015
$ws.Range("D7").Value = 0.001071841694489706
$ws.Range("E7").Value = 0.04427989729316284
$ws.Range("F7").Value = 0.6110633125826013
$ws.Range("G7").Value = 0.002363900624208227
$ws.Range("I7").Value = 0.4550797303182819
$ws.Range("M7").Value = 0.7824523813474968
$ws.Range("N7").Value = 1.513524774370069
$ws.Range("O7").Value = 2.001157253354307
$ws.Range("C8").Value = 0.004053348437640381
$ws.Range("D8").Value = 0.001033043347682838
$ws.Range("E8").Value = 0.04088257785811678
$ws.Range("F8").Value = 0.6584654399313621
$ws.Range("G8").Value = 0.002358491086496358
$ws.Range("I8").Value = 0.4904325128137259
$ws.Range("M8").Value = 0.9290972223864031
$ws.Range("N8").Value = 1.698341374432857
$ws.Range("O8").Value = 2.149177681767128
$ws.Range("C9").Value = 0.005402836275905543
$ws.Range("D9").Value = 0.0009667276668405833
$ws.Range("E9").Value = 0.03517626313417122
$ws.Range("F9").Value = 0.7553871655779574
$ws.Range("G9").Value = 0.002348954268788187
$ws.Range("I9").Value = 0.5627406488661393
$ws.Range("M9").Value = 1.215592491278514
$ws.Range("N9").Value = 2.061952193651962
$ws.Range("O9").Value = 2.452800009781242
$ws.Range("C10").Value = 0.006392005683252933
$ws.Range("D10").Value = 0.0009242040806656249
$ws.Range("E10").Value = 0.03158200434637148
$ws.Range("F10").Value = 0.8290870700990354
$ws.Range("G10").Value = 0.002342593460931394
$ws.Range("I10").Value = 0.6177395527227105
$ws.Range("M10").Value = 1.425383674267522
$ws.Range("N10").Value = 2.329625275652347
$ws.Range("O10").Value = 2.684266137338398
$ws.Range("C11").Value = 0.006841480855648285
$ws.Range("D11").Value = 0.0009062600076257254
$ws.Range("E11").Value = 0.03008013302187251
$ws.Range("F11").Value = 0.8631724695705572
$ws.Range("G11").Value = 0.002339838483040247
$ws.Range("I11").Value = 0.6431795206393502
$ws.Range("M11").Value = 1.52066466949465
$ws.Range("N11").Value = 2.451472790079322
$ws.Range("O11").Value = 2.791445286610326
$ws.Range("C12").Value = 0.007011608805512992
$ws.Range("D12").Value = 0.0008996711759481357
$ws.Range("E12").Value = 0.02953084795402816
$ws.Range("F12").Value = 0.8761612656802953
$ws.Range("G12").Value = 0.002338815055949817
$ws.Range("I12").Value = 0.652874370004497
$ws.Range("M12").Value = 1.55672191091071
$ws.Range("N12").Value = 2.497621287847664
$ws.Range("O12").Value = 2.832306230442498
$ws.Range("C13").Value = 0.006974972299794047
$ws.Range("D13").Value = 0.0009010809490765936
$ws.Range("E13").Value = 0.02964827701770734
$ws.Range("F13").Value = 0.8733602596506387
$ws.Range("G13").Value = 0.002339034589185766
$ws.Range("I13").Value = 0.6507836723789922
$ws.Range("M13").Value = 1.54895741140632
$ws.Range("N13").Value = 2.487682117247516
$ws.Range("O13").Value = 2.823493827813707
$ws.Range("C14").Value = 0.006855478994907571
$ws.Range("D14").Value = 0.0009057137792973791
$ws.Range("E14").Value = 0.0300345517070264
$ws.Range("F14").Value = 0.8642394284973847
$ws.Range("G14").Value = 0.002339753888134975
$ws.Range("I14").Value = 0.6439758889253824
$ws.Range("M14").Value = 1.523631601708473
$ws.Range("N14").Value = 2.455269328857923
$ws.Range("O14").Value = 2.794801418409236
$ws.Range("C15").Value = 0.006782275480631483
$ws.Range("D15").Value = 0.0009085785224804965
$ws.Range("E15").Value = 0.03027369657577861
$ws.Range("F15").Value = 0.8586632864924297
$ws.Range("G15").Value = 0.00234019705798949
$ws.Range("I15").Value = 0.6398139297587875
$ws.Range("M15").Value = 1.508115695538976
$ws.Range("N15").Value = 2.435416409958748
$ws.Range("O15").Value = 2.777262347948351
$ws.Range("C16").Value = 0.006362620770509864
$ws.Range("D16").Value = 0.000925405354690767
$ws.Range("E16").Value = 0.0316828597655201
$ws.Range("F16").Value = 0.8268708434465282
$ws.Range("G16").Value = 0.002342776284765506
$ws.Range("I16").Value = 0.6160855238314156
$ws.Range("M16").Value = 1.419153612496928
$ws.Range("N16").Value = 2.32166353239478
$ws.Range("O16").Value = 2.677299947049619
$ws.Range("C17").Value = 0.006105042950878214
$ws.Range("D17").Value = 0.0009360900881132395
$ws.Range("E17").Value = 0.03258165089374998
$ws.Range("F17").Value = 0.8075110913519836
$ws.Range("G17").Value = 0.002344393977042492
$ws.Range("I17").Value = 0.6016372251868347
$ws.Range("M17").Value = 1.364537756972595
$ws.Range("N17").Value = 2.251897699397148
$ws.Range("O17").Value = 2.616461356267735
$ws.Range("C18").Value = 0.005956844031928199
$ws.Range("D18").Value = 0.0009423671491943786
$ws.Range("E18").Value = 0.03311112929237847
$ws.Range("F18").Value = 0.7964284229676508
$ws.Range("G18").Value = 0.002345337481392867
$ws.Range("I18").Value = 0.5933664957234583
$ws.Range("M18").Value = 1.333109744997415
$ws.Range("N18").Value = 2.211778176876123
$ws.Range("O18").Value = 2.581645726420845
$ws.Range("C19").Value = 0.005906658547381483
$ws.Range("D19").Value = 0.0009445148916924495
$ws.Range("E19").Value = 0.03329254266702719
$ws.Range("F19").Value = 0.7926850162640022
$ws.Range("G19").Value = 0.00234565918050299
$ws.Range("I19").Value = 0.5905729375453177
$ws.Range("M19").Value = 1.322466328691235
$ws.Range("N19").Value = 2.198195890239788
$ws.Range("O19").Value = 2.56988806328917
$ws.Range("C20").Value = 0.006132467455309154
$ws.Range("D20").Value = 0.000934939036344673
$ws.Range("E20").Value = 0.03248467543763578
$ws.Range("F20").Value = 0.8095665240385728
$ws.Range("G20").Value = 0.002344220421159647
$ws.Range("I20").Value = 0.6031711722874178
$ws.Range("M20").Value = 1.37035321063577
$ws.Range("N20").Value = 2.259323600764674
$ws.Range("O20").Value = 2.622919371050045
$ws.Range("C21").Value = 0.006890579255603768
$ws.Range("D21").Value = 0.0009043473692509352
$ws.Range("E21").Value = 0.02992056351308126
$ws.Range("F21").Value = 0.8669162204955683
$ws.Range("G21").Value = 0.002339542075402129
$ws.Range("I21").Value = 0.6459738303423705
$ws.Range("M21").Value = 1.531071057941958
$ws.Range("N21").Value = 2.464789580901879
$ws.Range("O21").Value = 2.803221595161745
$ws.Range("C22").Value = 0.007385588461140458
$ws.Range("D22").Value = 0.0008855581981563887
$ws.Range("E22").Value = 0.02835822220647
$ws.Range("F22").Value = 0.9048724738977967
$ws.Range("G22").Value = 0.002336600004583114
$ws.Range("I22").Value = 0.6743054325277313
$ws.Range("M22").Value = 1.635971227548396
$ws.Range("N22").Value = 2.599115192116528
$ws.Range("O22").Value = 2.92266136087926
$ws.Range("C23").Value = 0.007121436934248493
$ws.Range("D23").Value = 0.0008954744969029704
$ws.Range("E23").Value = 0.02918159532040043
$ws.Range("F23").Value = 0.8845707137802918
$ws.Range("G23").Value = 0.002338159709512006
$ws.Range("I23").Value = 0.659151340825062
$ws.Range("M23").Value = 1.579997153483447
$ws.Range("N23").Value = 2.527420659581423
$ws.Range("O23").Value = 2.858766344333162
$ws.Range("C24").Value = 0.006120069191489108
$ws.Range("D24").Value = 0.0009354590083505698
$ws.Range("E24").Value = 0.03252847834317407
$ws.Range("F24").Value = 0.8086371150076133
$ws.Range("G24").Value = 0.002344298843830028
$ws.Range("I24").Value = 0.6024775633824078
$ws.Range("M24").Value = 1.367724133104119
$ws.Range("N24").Value = 2.255966382249369
$ws.Range("O24").Value = 2.619999200452185
$ws.Range("C25").Value = 0.005038157547222255
$ws.Range("D25").Value = 0.0009835974361682709
$ws.Range("E25").Value = 0.03661605557468661
$ws.Range("F25").Value = 0.7287350745715457
$ws.Range("G25").Value = 0.00235142028869243
$ws.Range("I25").Value = 0.5428543734047935
$ws.Range("M25").Value = 1.138207394449481
$ws.Range("N25").Value = 2.061952193651962
$ws.Range("O25").Value = 2.072273628073162
